# Updated the three sample schemas with preparation_condition and storage_method
#
# Renames:
#   "preparation_temperature list" -> "preparation_condition list"
#   "storage_temperature list"     -> "storage_method list"
# and updates their contents, the "Export as TSV" sheet's column G / K
# comments and data validation rules to match.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the two list sheets
# ------------------------------------------------------------------
$wsPrep = $wb.Worksheets.Item("preparation_temperature list")
$wsPrep.Name = "preparation_condition list"

$wsStorage = $wb.Worksheets.Item("storage_temperature list")
$wsStorage.Name = "storage_method list"

# ------------------------------------------------------------------
# 2. Replace the contents of the preparation_condition list sheet
# ------------------------------------------------------------------
$prepValues = @(
    "frozen in liquid nitrogen",
    "frozen in liquid nitrogen vapor",
    "frozen in ice",
    "frozen in dry ice",
    "frozen at -20 C",
    "ambient temperature",
    "unknown"
)
for ($i = 0; $i -lt $prepValues.Length; $i++) {
    $wsPrep.Cells.Item($i + 1, 1).Value = $prepValues[$i]
}
# the old list had 8 rows, the new one only has 7 - drop the leftover row
$wsPrep.Rows("8").Delete()

# ------------------------------------------------------------------
# 3. Replace the contents of the storage_method list sheet
# ------------------------------------------------------------------
$storageValues = @(
    "frozen in liquid nitrogen",
    "frozen in liquid nitrogen vapor",
    "frozen in ice",
    "frozen in dry ice",
    "frozen at -80 C",
    "frozen at -20 C",
    "refrigerator",
    "ambient temperature",
    "incubated at 37 C",
    "none",
    "unknown"
)
for ($i = 0; $i -lt $storageValues.Length; $i++) {
    $wsStorage.Cells.Item($i + 1, 1).Value = $storageValues[$i]
}
# the old list had 12 rows, the new one only has 11 - drop the leftover row
$wsStorage.Rows("12").Delete()

# ------------------------------------------------------------------
# 4. Update the "Export as TSV" sheet: column header names and comments
# ------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Export as TSV")

$wsMain.Range("G1").Value = "preparation_condition"
$wsMain.Range("K1").Value = "storage_method"

$wsMain.Range("G1").Comment.Text("The condition under which the preparation occurred, such as whether the sample was placed in dry ice during the preparation.")
$wsMain.Range("K1").Comment.Text("The method by which the sample was stored, after preparation and before the assay was performed.")

# ------------------------------------------------------------------
# 5. Update the data validation rules that reference the two lists
# ------------------------------------------------------------------
$gValidation = $wsMain.Range("G2:G1048576").Validation
$gValidation.Formula1 = "'preparation_condition list'!`$A`$1:`$A`$7"
$gValidation.ErrorMessage = "Value must come from preparation_condition list."

$kValidation = $wsMain.Range("K2:K1048576").Validation
$kValidation.Formula1 = "'storage_method list'!`$A`$1:`$A`$11"
$kValidation.ErrorMessage = "Value must come from storage_method list."
